$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 157; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = $cell.Value2
    $parts = $orig.Split(",")
    $n = $parts.Length
    if ($n -gt 1) {
        $result = $parts[1].Trim()
        for ($i = 2; $i -lt $n; $i++) {
            $result = $result + ", " + $parts[$i].Trim()
        }
        $result = $result + ", " + $parts[0].Trim()
        $cell.Value = $result
    }
}
